$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 for "Business demography"; this pushes the
# existing Key Stage 4 / Key Stage 5 rows down from 11/12 to 12/13.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the Business demography data.
$ws.Range("A11").Value = "Business demography"
$ws.Range("B11").Value = "<a href=' https://www.ons.gov.uk/businessindustryandtrade/business/activitysizeandlocation/datasets/businessdemographyreferencetable'>ONS Business Demography</a>"
$ws.Range("C11").Value = "2021 (17/11/22)"
$ws.Range("D11").Value = "2022 (16/11/23)"

# Correct "ONS UK Business Count" -> "ONS UK Business Counts" for the
# Enterprise by employment size rows (now rows 9 and 10).
$ws.Range("B9").Value = "<a href='https://www.nomisweb.co.uk/datasets/idbrent'>ONS UK Business Counts</a>"
$ws.Range("B10").Value = "<a href='https://www.nomisweb.co.uk/datasets/idbrent'>ONS UK Business Counts</a>"

# Update the sheet view to match the authored state: scrolled to column B,
# with C12 selected (instead of the old A2:A12 selection anchored at A12).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C12").Select()
